$d = $word.ActiveDocument

# The section has a "primary"/default header & footer plus a "first page"
# header & footer (wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2).
# Each one holds exactly one inline picture:
#   - Header pictures (the BTEC logo) are currently named "image1.jpg" and
#     need to become "image2.jpg".
#   - Footer pictures (the Pearson Edexcel logo) are currently named
#     "image2.png" and need to become "image1.png".
# Renaming is done by selecting the picture's range and renaming it through
# $word.Selection.InlineShapes - going through the Selection keeps the
# rename reliable for header *and* footer stories alike.

function Rename-InlinePicture($headerFooter, $newName) {
    if (-not $headerFooter.Exists) {
        return
    }
    $shapes = $headerFooter.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shapes.Item($i).Range.Select()
        $word.Selection.InlineShapes.Item(1).Name = $newName
    }
}

foreach ($idx in 1, 2) {
    $sec = $d.Sections.Item(1)
    Rename-InlinePicture $sec.Headers.Item($idx) "image2.jpg"

    $sec = $d.Sections.Item(1)
    Rename-InlinePicture $sec.Footers.Item($idx) "image1.png"
}
